$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the 1.000.000-records block (old rows 9-11
# shift down to 11-13), making room for new SHA3_512 / SHA3_384 rows.
$ws.Rows("9:10").Insert()

# Insert two new rows before the 2.000.000-records block (old rows 13-15,
# already shifted to 15-17 by the previous insert, shift down to 17-19),
# making room for new SHA3_512 / SHA3_384 rows.
$ws.Rows("15:16").Insert()

# --- New row 9: SHA3_512, 1.000.000 records ---
$ws.Range("A9").Value = "SHA3_512"
$ws.Range("B9").Value = 1000000
$ws.Range("C9").Value = 14
$ws.Range("D9").Value = 9828
$ws.Range("E9").Value = 10317
$ws.Range("F9").Value = 9784
$ws.Range("G9").Formula = "=AVERAGE(D9:F9)"
$ws.Range("H9").Value = 11752
$ws.Range("I9").Value = 11436
$ws.Range("J9").Value = 11257
$ws.Range("K9").Formula = "=AVERAGE(H9:J9)"
$ws.Range("L9").Formula = "=(K9*500)/1000/60/60"

# --- New row 10: SHA3_384, 1.000.000 records (no measurements, errors) ---
$ws.Range("A10").Value = "SHA3_384"
$ws.Range("B10").Value = 1000000
$ws.Range("C10").Value = 14
$ws.Range("G10").Formula = "=AVERAGE(D10:F10)"
$ws.Range("K10").Formula = "=AVERAGE(H10:J10)"
$ws.Range("L10").Formula = "=(K10*500)/1000/60/60"

# --- New row 15: SHA3_512, 2.000.000 records ---
$ws.Range("A15").Value = "SHA3_512"
$ws.Range("B15").Value = 2000000
$ws.Range("C15").Value = 14
$ws.Range("D15").Value = 19946
$ws.Range("E15").Value = 19402
$ws.Range("F15").Value = 20340
$ws.Range("G15").Formula = "=AVERAGE(D15:F15)"
$ws.Range("H15").Value = 21951
$ws.Range("I15").Value = 24499
$ws.Range("J15").Value = 24550
$ws.Range("K15").Formula = "=AVERAGE(H15:J15)"
$ws.Range("L15").NumberFormat = "0.000"
$ws.Range("L15").Formula = "=(K15*250)/1000/60/60"

# --- New row 16: SHA3_384, 2.000.000 records (no measurements, errors) ---
$ws.Range("A16").Value = "SHA3_384"
$ws.Range("B16").Value = 2000000
$ws.Range("C16").Value = 14
$ws.Range("G16").Formula = "=AVERAGE(D16:F16)"
$ws.Range("K16").Formula = "=AVERAGE(H16:J16)"
$ws.Range("L16").NumberFormat = "0.000"
$ws.Range("L16").Formula = "=(K16*250)/1000/60/60"

# Restore the selection state as saved by Excel after the edit.
$ws.Range("L20").Select()
